# [Fonds de solidarite] Add 2020-12-09 data
# Update the "nombre_aides" (C) and "montant_total" (D) columns for the
# rows whose underlying figures changed in the 2020-12-09 refresh.
# Values are stored as plain text in the source data (e.g. "1864524.90"),
# so we force text formatting before writing to avoid Excel silently
# normalising the decimal representation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 8;   C = "771";  D = "1866524.90" },
    @{ Row = 9;   C = "1364"; D = "10197813.05" },
    @{ Row = 62;  C = "75";   D = "714460.00" },
    @{ Row = 127; C = "337";  D = "1115116.00" },
    @{ Row = 129; C = "1120"; D = "8961082.69" },
    @{ Row = 133; C = "162";  D = "1442979.82" },
    @{ Row = 145; C = "8228"; D = "25800767.91" },
    @{ Row = 146; C = "5005"; D = "32769305.85" },
    @{ Row = 151; C = "849";  D = "3233731.22" },
    @{ Row = 244; C = "475";  D = "3464962.16" },
    @{ Row = 259; C = "604";  D = "1513025.18" }
)

foreach ($u in $updates) {
    $cCell = $ws.Cells.Item($u.Row, 3)
    $cCell.NumberFormat = "@"
    $cCell.Value = $u.C

    $dCell = $ws.Cells.Item($u.Row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $u.D
}
